$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

$ws.Range("A2").Value = "mngr506568"
$ws.Range("B2").Value = "sYdAjun"
